$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the indicator reporting details.
# Set values in the same order they were entered/appended so that the
# shared-strings table is rebuilt in the same order as the target workbook.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B7").Value = "Mambetaliev T.A."
$ws.Range("B9").Value = "(0312) 62 56 07"
$ws.Range("B4").Value = "11.5.2 Direct economic loss attributed to disasters in relation to global gross domestic product (GDP)"

# Move the active selection to B4 (also drops the stale topLeftCell scroll anchor).
$ws.Range("B4").Select()
